{"js": "// The document has a \"Requisitos\" section ending with a paragraph that\n// reads \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\". Immediately\n// after it there used to be four extra paragraphs (an empty spacer, a\n// \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, another empty\n// spacer, and an empty page-break paragraph) that must be removed,\n// leaving only the original trailing empty paragraph and the trailing\n// page-break paragraph that close the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"LOB1036: Geometria Anal\u00edtica ...\") by\n// its text instead of a hard-coded index, so the script is resilient to\n// unrelated edits earlier in the document.\nconst anchorText = \"LOB1036: Geometria Anal\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOB1036' anchor paragraph.\");\n}\n\n// The four paragraphs right after the anchor are the ones to drop:\n//   1) empty spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) empty spacer paragraph\n//   4) empty paragraph with pageBreakBefore + jc=left\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 4 && i < paragraphs.items.length; i++) {\n  toDelete.push(paragraphs.items[i]);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The document has a \"Requisitos\" section ending with a paragraph that\n# reads \"LOB1036: Geometria Anal\u00edtica (Requisito fraco)\". Immediately\n# after it there used to be four extra paragraphs (an empty spacer, a\n# \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, another empty\n# spacer, and an empty page-break paragraph) that must be removed,\n# leaving only the original trailing empty paragraph and the trailing\n# page-break paragraph that close the document.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"LOB1036: Geometria Anal\u00edtica ...\") by\n# its text instead of a hard-coded index, so the script is resilient to\n# unrelated edits earlier in the document.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*LOB1036: Geometria Anal*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOB1036' anchor paragraph.\"\n}\n\n# The four paragraphs right after the anchor are the ones to drop:\n#   1) empty spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) empty spacer paragraph\n#   4) empty paragraph with pageBreakBefore + jc=left\n# Delete from the end backwards so earlier indices stay valid.\nfor ($j = $anchorIndex + 4; $j -ge $anchorIndex + 1; $j--) {\n    $d.Paragraphs.Item($j).Range.Delete()\n}\n"}
